# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet, with
#    the same column layout as the other quarterly sheets (基金代码 /
#    基金名称 / 基金规模 / 股票总仓位 / 仓位占比 / 持有市值(亿元) / 仓位排名)
#    and a single holding row (鹏华中证传媒指数（LOF）).
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing quarters down by one row and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet just before "总计"
# ---------------------------------------------------------------------
$zjBeforeAdd = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($zjBeforeAdd)
$q1.Name = "2022-Q1"

# NOTE: after Worksheets.Add(before), the handle used as the "before"
# argument can end up aliased to the newly created sheet, so re-resolve
# "总计" by name to get a fresh, correct reference to it.
$zj = $wb.Worksheets.Item("总计")

# Reuse the exact cell formatting (fonts/borders/alignment) of the
# previous quarter sheet so the new sheet matches the established style.
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$prevQuarter.Range("B1:H2").Copy()
$q1.Range("B1").PasteSpecial(-4122)
$prevQuarter.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

function Set-TextValue($range, $value) {
    # Force the written value to be stored as text (matching the source
    # data's convention of keeping numeric-looking fund figures as
    # strings) without leaving a lingering custom number format behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "160629"
Set-TextValue $q1.Range("C2") "鹏华中证传媒指数（LOF）"
Set-TextValue $q1.Range("D2") "7.63"
Set-TextValue $q1.Range("E2") "92.90"
Set-TextValue $q1.Range("F2") "3.17"
Set-TextValue $q1.Range("G2") "0.2419"
$q1.Range("H2").Value = 9

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q1" row at the top of "总计"
# ---------------------------------------------------------------------
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()

# Copy the index cell's formatting (style s="2") down into the freshly
# inserted A2 cell.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.24

# Renumber the index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3
$zj.Range("A6").Value = 4
$zj.Range("A7").Value = 5

Write-Host "2022-Q1 sheet added and summary sheet updated"
